$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. The "last saved" date field (shown on every slide layout and on the
#    slide master) moved on from 09/05/2022 to 18/10/2023.
# ---------------------------------------------------------------------
$newDate = "18/10/2023"

$sm = $p.SlideMaster
for ($j = 1; $j -le $sm.Shapes.Count; $j++) {
    $sh = $sm.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $sm.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2. "Corrected gating scheme based on what was labelled" - on the
#    gating-scheme slide, the little "CD3" caption was actually sitting
#    on top of what is labelled "Plot 1", so the stray "CD3" textbox is
#    removed and the highlighted gate rectangle is widened to cover the
#    area that is really being called out.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(2)

$cd3Box = $null
$gateRect = $null
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $sh = $s.Shapes.Item($j)
    if ($sh.Id -eq 12) { $cd3Box = $sh }
    if ($sh.Id -eq 46) { $gateRect = $sh }
}

if ($cd3Box -ne $null) {
    $cd3Box.Delete()
}

if ($gateRect -ne $null) {
    $gateRect.Left = 265.622635
    $gateRect.Width = 139.710355
}
